$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 57214.285
$ws.Range("J3").Value = 57214.285
$ws.Range("L3").Value = 57214.285
$ws.Range("N3").Value = -57442.285

$ws.Range("H6").Value = 2220
$ws.Range("I6").Value = 2700
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 8100
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = -7988
$ws.Range("N6").Value = -1124

$ws.Range("H95").Value = 40375
$ws.Range("J95").Value = 40375
$ws.Range("L95").Value = 40375
$ws.Range("N95").Value = -45867

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H98").Value = 21354.777
$ws.Range("I98").Value = 19024.125
$ws.Range("K98").Value = 19024.125
$ws.Range("M98").Value = -17526.125

$ws.Range("H102").Value = 57214.285
$ws.Range("J102").Value = 57214.285
$ws.Range("L102").Value = 57214.285
$ws.Range("N102").Value = -63704.285

$ws.Range("H122").Value = 21354.777
$ws.Range("I122").Value = 19024.125
$ws.Range("K122").Value = 57072.375
$ws.Range("M122").Value = -54622.375

$ws.Range("H125").Value = 2507.75
$ws.Range("I125").Value = 2016
$ws.Range("J125").Value = 2999.5
$ws.Range("K125").Value = 18144
$ws.Range("L125").Value = 26995.5
$ws.Range("M125").Value = -15684
$ws.Range("N125").Value = -31915.5

$ws.Range("H135").Value = 250
$ws.Range("I135").Value = 250
$ws.Range("K135").Value = 2250
$ws.Range("M135").Value = 285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 119999.5
$ws.Range("J92").Value = 119999.5
$ws.Range("L92").Value = 119999.5
$ws.Range("N92").Value = -124991.5

$ws.Range("H122").Value = 9488.700000000001
$ws.Range("J122").Value = 8278.6
$ws.Range("L122").Value = 24835.8
$ws.Range("N122").Value = -29735.8

$ws.Range("H132").Value = 1920
$ws.Range("I132").Value = 1920
$ws.Range("K132").Value = 5760
$ws.Range("M132").Value = -3230

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 70000
$ws.Range("J68").Value = 70000
$ws.Range("L68").Value = 70000
$ws.Range("N68").Value = -71622

$ws.Range("H71").Value = 70000
$ws.Range("J71").Value = 70000
$ws.Range("L71").Value = 210000
$ws.Range("N71").Value = -218112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 42999.75
$ws.Range("I56").Value = 42999.75
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 42999.75
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -42154.75
$ws.Range("N56").ClearContents()

$ws.Range("H58").Value = 842
$ws.Range("I58").Value = 842
$ws.Range("K58").Value = 842
$ws.Range("M58").Value = -639

$ws.Range("H93").Value = 19999
$ws.Range("I93").Value = 19999
$ws.Range("K93").Value = 19999
$ws.Range("M93").Value = -18127

$ws.Range("H95").Value = 36705.285
$ws.Range("J95").Value = 36705.285
$ws.Range("L95").Value = 36705.285
$ws.Range("N95").Value = -42197.285

$ws.Range("H136").Value = 842
$ws.Range("I136").Value = 842
$ws.Range("K136").Value = 2526
$ws.Range("M136").Value = 24

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 11150.5
$ws.Range("J34").Value = 11700.533
$ws.Range("L34").Value = 35101.599
$ws.Range("N34").Value = -35269.599

$ws.Range("H39").Value = 10034.333
$ws.Range("I39").Value = 51.5
$ws.Range("K39").Value = 154.5
$ws.Range("M39").Value = 139.5

$ws.Range("H55").Value = 1758.721
$ws.Range("I55").Value = 1675
$ws.Range("J55").Value = 1767.3077
$ws.Range("K55").Value = 5025
$ws.Range("L55").Value = 5301.9231
$ws.Range("M55").Value = -4848
$ws.Range("N55").Value = -5655.9231

$ws.Range("H129").Value = 1210.091
$ws.Range("I129").Value = 295.4
$ws.Range("J129").Value = 1972.3334
$ws.Range("K129").Value = 886.1999999999999
$ws.Range("L129").Value = 5917.0002
$ws.Range("M129").Value = 4113.8
$ws.Range("N129").Value = -15917.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H122").Value = 6054.8887
$ws.Range("I122").Value = 5666.6665
$ws.Range("J122").Value = 6249
$ws.Range("K122").Value = 16999.9995
$ws.Range("L122").Value = 18747
$ws.Range("M122").Value = -14549.9995
$ws.Range("N122").Value = -23647

$ws.Range("H132").Value = 4142.2
$ws.Range("I132").Value = 2903.6667
$ws.Range("K132").Value = 8711.000100000001
$ws.Range("M132").Value = -6181.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 26663.334
$ws.Range("I50").Value = 42990
$ws.Range("J50").Value = 18500
$ws.Range("K50").Value = 42990
$ws.Range("L50").Value = 18500
$ws.Range("M50").Value = -42353
$ws.Range("N50").Value = -19774

$ws.Range("H63").Value = 90000
$ws.Range("I63").Value = 90000
$ws.Range("K63").Value = 90000
$ws.Range("M63").Value = -89251

$ws.Range("H66").Value = 90000
$ws.Range("I66").Value = 90000
$ws.Range("K66").Value = 270000
$ws.Range("M66").Value = -266256

$ws.Range("H74").Value = 100060.664
$ws.Range("I74").Value = 100060.664
$ws.Range("K74").Value = 100060.664
$ws.Range("M74").Value = -99062.664

$ws.Range("H77").Value = 100060.664
$ws.Range("I77").Value = 100060.664
$ws.Range("K77").Value = 300181.992
$ws.Range("M77").Value = -295189.992

$ws.Range("H82").Value = 4250
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 4250
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws.Range("H95").Value = 14565.333
$ws.Range("J95").Value = 14565.333
$ws.Range("L95").Value = 14565.333
$ws.Range("N95").Value = -20057.333

$ws.Range("H132").Value = 2820.2
$ws.Range("I132").Value = 3583.8333
$ws.Range("J132").Value = 2311.111
$ws.Range("K132").Value = 10751.4999
$ws.Range("L132").Value = 6933.333
$ws.Range("M132").Value = -8221.499899999999
$ws.Range("N132").Value = -11993.333

$ws.Range("H136").Value = 1264499
$ws.Range("I136").Value = 1264499
$ws.Range("K136").Value = 3793497
$ws.Range("M136").Value = -3790947

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2388.4285
$ws.Range("I23").Value = 273.33334
$ws.Range("J23").Value = 3974.75
$ws.Range("K23").Value = 273.33334
$ws.Range("L23").Value = 3974.75
$ws.Range("M23").Value = -44.33334000000002
$ws.Range("N23").Value = -4432.75

$ws.Range("H122").Value = 3499
$ws.Range("I122").Value = 2999.5
$ws.Range("K122").Value = 8998.5
$ws.Range("M122").Value = -6548.5
